$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new student records below the existing table (rows 56-58) ---
# Copy the formatting of the last existing data row (55) down into the three
# new rows so the new records look consistent with the rest of the table.
$ws.Range("A55:G55").Copy()
$ws.Range("A56:G58").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column A (Student_ID) holds numeric-looking values but is stored as TEXT
# in this workbook, so force a text number format before writing the value.
$ws.Range("A56:A58").NumberFormat = "@"

# Row 56 - Severus Snape
$ws.Range("A56").Value = "242743452"
$ws.Range("B56").Value = "Severus"
$ws.Range("C56").ClearContents()
$ws.Range("D56").Value = "Snape"
$ws.Range("E56").Value = 20020414
$ws.Range("F56").Value = "Xenobotany Society"
$ws.Range("G56").Value = 72.5

# Row 57 - Walter Hartwell White
$ws.Range("A57").Value = "224303042"
$ws.Range("B57").Value = "Walter"
$ws.Range("C57").Value = "Hartwell"
$ws.Range("D57").Value = "White"
$ws.Range("E57").Value = 20060510
$ws.Range("F57").Value = "ABW"
$ws.Range("G57").Value = 71.5

# Row 58 - Frodo M Baggins
$ws.Range("A58").Value = "212954131"
$ws.Range("B58").Value = "Frodo"
$ws.Range("C58").Value = "M"
$ws.Range("D58").Value = "Baggins"
$ws.Range("E58").Value = 20040714
$ws.Range("F58").Value = "The Temporal Anomaly Watch"
$ws.Range("G58").Value = 33

# --- Update the view state to match where the user ended up working ---
$ws.Range("C62").Select()

Write-Output "done"
